$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: operations are applied from the bottom of the sheet upward so that
# earlier (lower) row numbers referenced below remain valid while we still
# need them.

# Account numbers ("Conta" column) are zero-padded numeric strings, so the
# target cells must be formatted as Text before assignment; otherwise Excel
# auto-converts them to numbers and the leading zeros are lost.

# 7. Insert a new row for "004214592 / MERG / 2516.26" right after the
#    ERIKA row (004971783, row 15) i.e. before row 16 (001761119 / BLUEMETRIX)
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "004214592"
$ws.Cells.Item(16, 2).Value = "MERG"
$ws.Cells.Item(16, 3).Value = 2516.26

# 6. Insert two new rows before row 12 (004224011 / THOMAS):
#    "005624730 / ISABEL / 20000" and "004480970 / ALBERTO / 16352.97"
$ws.Rows.Item(12).Resize(2).Insert()
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "005624730"
$ws.Cells.Item(12, 2).Value = "ISABEL"
$ws.Cells.Item(12, 3).Value = 20000
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "004480970"
$ws.Cells.Item(13, 2).Value = "ALBERTO"
$ws.Cells.Item(13, 3).Value = 16352.97

# 5. Delete row 10 (004895776 / FERNANDO) entirely
$ws.Rows.Item(10).Delete()

# 4. Insert a new row before row 9 (004321016 / JOAQUIM):
#    "004444380 / MARCELO / 51364.77"
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "004444380"
$ws.Cells.Item(9, 2).Value = "MARCELO"
$ws.Cells.Item(9, 3).Value = 51364.77

# 3. Update row 7 (005277762 / NIVALDO) balance from 114450 to 64450
$ws.Cells.Item(7, 3).Value = 64450

# 2. Update row 3 (005103059 / WALQUIRIA) balance from 228900 to 239243.88
$ws.Cells.Item(3, 3).Value = 239243.88

# 1. Insert a new row before row 2 (005642649 / VR):
#    "005320069 / RICARDO / 1114376.75"
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "005320069"
$ws.Cells.Item(2, 2).Value = "RICARDO"
$ws.Cells.Item(2, 3).Value = 1114376.75
